$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2

# 1) High-income / zero-car-owning paragraph: reorder clauses and change
#    "there is" -> "have".
$old1 = "In high-income earning, most walkable communities and the percentage of zero car owning households on a census block level there is a weak and negative correlation. "
$new1 = "In high-income earning and most walkable communities, the percentage of zero car owning households on a census block level have a weak and negative correlation. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) High-income / one-car-owning paragraph: reorder clauses, keep "have" and
#    "is" both present (as in the authored edit).
$old2 = "In high-income earning, most walkable communities and the percentage of one car owning households on a census block level there is a weak and positive correlation. "
$new2 = "In high-income earning and most walkable communities, the percentage of one car owning households on a census block level have is a weak and positive correlation. "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3) High-income / two-plus-car-owning paragraph: reorder clauses, keep
#    "have is" both present (as in the authored edit).
$old3 = "In high-income earning, most walkable communities and percentage of households with two plus cars ownership on a census block level there is a weak and negative correlation. "
$new3 = "In high-income earning and most walkable communities, percentage of households with two plus cars ownership on a census block level have is a weak and negative correlation. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# 4) Transit-stop sample-size paragraph: the visible text does not change, but
#    the grammar-check markers around "the" and "portrays" are cleared, and
#    the runs they split apart are merged back into one run. We rebuild just
#    that stretch of text (leaving the preceding "centralized transit" run
#    untouched) so only the intended run boundary goes away.
$old4 = " stops in meters and car ownership in the most walkable communities. We used a random sample size of 2000 since the dataset is dense and charted results in a scatter plot. This sample size portrays that individuals in walkable communities will own a car whether they are near a transit stop or not."

$findRng = $d.Content
$findRng.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$segStart = $findRng.Start
$segEnd = $findRng.End

$segment = $d.Range($segStart, $segEnd)
$segment.Text = ""

$insertPoint = $d.Range($segStart, $segStart)
$insertPoint.InsertAfter($old4)

$newRun = $d.Range($segStart, $segStart + $old4.Length)
$newRun.Font.Name = "Calibri"
$newRun.Font.NameAscii = "Calibri"
$newRun.Font.NameOther = "Calibri"
$newRun.Font.Color = 0
